$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3:F38").Formula = "=E3-E2"

$ws.Range("F1").Value = "d_un"
$ws.Range("H1").Value = "d_inf"
$ws.Range("I1").Value = "cpi94"
$ws.Range("J1").Value = "ngdp"

[void]$ws.Range("O11").Select()
